$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Jisun"
$ws.Range("C6").Value = 160
$ws.Range("D6").Value = "Yellow"
$ws.Activate() | Out-Null
$ws.Range("D6").Select() | Out-Null
